{"js": "// Replace the 25 division problems in the table with their new values.\n// Each \"before\" text is unique within the document, so an exact\n// case-sensitive search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"402\u00f77=\", \"210\u00f79=\"],\n  [\"897\u00f72=\", \"799\u00f77=\"],\n  [\"673\u00f74=\", \"621\u00f75=\"],\n  [\"279\u00f74=\", \"664\u00f73=\"],\n  [\"952\u00f76=\", \"795\u00f75=\"],\n  [\"776\u00f74=\", \"120\u00f73=\"],\n  [\"501\u00f77=\", \"802\u00f77=\"],\n  [\"308\u00f72=\", \"502\u00f75=\"],\n  [\"182\u00f78=\", \"950\u00f76=\"],\n  [\"623\u00f77=\", \"557\u00f75=\"],\n  [\"123\u00f76=\", \"135\u00f79=\"],\n  [\"770\u00f72=\", \"932\u00f79=\"],\n  [\"638\u00f77=\", \"881\u00f74=\"],\n  [\"876\u00f77=\", \"591\u00f74=\"],\n  [\"778\u00f79=\", \"394\u00f77=\"],\n  [\"435\u00f74=\", \"271\u00f75=\"],\n  [\"250\u00f72=\", \"853\u00f76=\"],\n  [\"430\u00f79=\", \"652\u00f78=\"],\n  [\"338\u00f79=\", \"167\u00f77=\"],\n  [\"115\u00f76=\", \"548\u00f74=\"],\n  [\"649\u00f74=\", \"327\u00f72=\"],\n  [\"259\u00f77=\", \"280\u00f74=\"],\n  [\"210\u00f74=\", \"949\u00f76=\"],\n  [\"165\u00f74=\", \"116\u00f79=\"],\n  [\"405\u00f77=\", \"357\u00f75=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"402\u00f77=\", \"210\u00f79=\"),\n    @(\"897\u00f72=\", \"799\u00f77=\"),\n    @(\"673\u00f74=\", \"621\u00f75=\"),\n    @(\"279\u00f74=\", \"664\u00f73=\"),\n    @(\"952\u00f76=\", \"795\u00f75=\"),\n    @(\"776\u00f74=\", \"120\u00f73=\"),\n    @(\"501\u00f77=\", \"802\u00f77=\"),\n    @(\"308\u00f72=\", \"502\u00f75=\"),\n    @(\"182\u00f78=\", \"950\u00f76=\"),\n    @(\"623\u00f77=\", \"557\u00f75=\"),\n    @(\"123\u00f76=\", \"135\u00f79=\"),\n    @(\"770\u00f72=\", \"932\u00f79=\"),\n    @(\"638\u00f77=\", \"881\u00f74=\"),\n    @(\"876\u00f77=\", \"591\u00f74=\"),\n    @(\"778\u00f79=\", \"394\u00f77=\"),\n    @(\"435\u00f74=\", \"271\u00f75=\"),\n    @(\"250\u00f72=\", \"853\u00f76=\"),\n    @(\"430\u00f79=\", \"652\u00f78=\"),\n    @(\"338\u00f79=\", \"167\u00f77=\"),\n    @(\"115\u00f76=\", \"548\u00f74=\"),\n    @(\"649\u00f74=\", \"327\u00f72=\"),\n    @(\"259\u00f77=\", \"280\u00f74=\"),\n    @(\"210\u00f74=\", \"949\u00f76=\"),\n    @(\"165\u00f74=\", \"116\u00f79=\"),\n    @(\"405\u00f77=\", \"357\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
